# Refresh crypto price/volume snapshot (GitHub Actions scheduled update).
# New values are prefixed with a leading apostrophe so Excel stores them
# as text (matching the source data, which includes thousands-separator
# "prices" like 67.139.21 that must not be reinterpreted as numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '67.139.21' }
    @{ Cell = 'E2'; Value = '  -0.69%  ' }
    @{ Cell = 'D3'; Value = '2.472.69' }
    @{ Cell = 'E3'; Value = '  -0.71%  ' }
    @{ Cell = 'E4'; Value = '  -0.08%  ' }
    @{ Cell = 'D5'; Value = '582.59' }
    @{ Cell = 'E5'; Value = '  -1.22%  ' }
    @{ Cell = 'D6'; Value = '168.95' }
    @{ Cell = 'E6'; Value = '  -2.33%  ' }
    @{ Cell = 'E7'; Value = '  -0.03%  ' }
    @{ Cell = 'E8'; Value = '  -1.61%  ' }
    @{ Cell = 'D9'; Value = '2.473.12' }
    @{ Cell = 'E9'; Value = '  -0.59%  ' }
    @{ Cell = 'E10'; Value = '  -2.94%  ' }
    @{ Cell = 'E11'; Value = '  -0.88%  ' }
    @{ Cell = 'D12'; Value = '4.96' }
    @{ Cell = 'E12'; Value = '  -2.65%  ' }
    @{ Cell = 'E13'; Value = '  -2.17%  ' }
    @{ Cell = 'B14'; Value = 'Avalanche' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax' }
    @{ Cell = 'D14'; Value = '25.56' }
    @{ Cell = 'E14'; Value = '  -2.60%  ' }
    @{ Cell = 'B15'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D15'; Value = '2.920.33' }
    @{ Cell = 'E15'; Value = '  -1.16%  ' }
    @{ Cell = 'D16'; Value = '67.117.52' }
    @{ Cell = 'E16'; Value = '  -0.69%  ' }
    @{ Cell = 'D18'; Value = '2.459.89' }
    @{ Cell = 'E18'; Value = '  -2.14%  ' }
    @{ Cell = 'D19'; Value = '11.26' }
    @{ Cell = 'E19'; Value = '  -4.28%  ' }
    @{ Cell = 'D20'; Value = '7.56' }
    @{ Cell = 'E20'; Value = '  -5.19%  ' }
    @{ Cell = 'D21'; Value = '355.81' }
    @{ Cell = 'E21'; Value = '  -3.05%  ' }
    @{ Cell = 'E22'; Value = '  -1.84%  ' }
    @{ Cell = 'E23'; Value = '  +0.40%  ' }
    @{ Cell = 'D24'; Value = '69.30' }
    @{ Cell = 'E24'; Value = '  -3.10%  ' }
    @{ Cell = 'E26'; Value = '  -6.24%  ' }
    @{ Cell = 'D27'; Value = '9.12' }
    @{ Cell = 'E27'; Value = '  -7.94%  ' }
    @{ Cell = 'D28'; Value = '0.998' }
    @{ Cell = 'E28'; Value = '  -0.53%  ' }
    @{ Cell = 'D29'; Value = '2.599.02' }
    @{ Cell = 'E29'; Value = '  -0.90%  ' }
    @{ Cell = 'E30'; Value = '  -5.46%  ' }
    @{ Cell = 'D31'; Value = '510.35' }
    @{ Cell = 'E31'; Value = '  -4.01%  ' }
    @{ Cell = 'E32'; Value = '  -7.22%  ' }
    @{ Cell = 'E33'; Value = '  -4.31%  ' }
    @{ Cell = 'E34'; Value = '  -5.52%  ' }
    @{ Cell = 'E35'; Value = '  -0.10%  ' }
    @{ Cell = 'E36'; Value = '  -6.30%  ' }
    @{ Cell = 'D37'; Value = '159.95' }
    @{ Cell = 'E37'; Value = '  +1.27%  ' }
    @{ Cell = 'E38'; Value = '  -0.23%  ' }
    @{ Cell = 'D39'; Value = '18.43' }
    @{ Cell = 'E39'; Value = '  -1.50%  ' }
    @{ Cell = 'E40'; Value = '  -5.26%  ' }
    @{ Cell = 'E41'; Value = '  +0.18%  ' }
    @{ Cell = 'E42'; Value = '  -6.04%  ' }
    @{ Cell = 'E43'; Value = '  -6.00%  ' }
    @{ Cell = 'D44'; Value = '0.326' }
    @{ Cell = 'E44'; Value = '  -6.37%  ' }
    @{ Cell = 'E45'; Value = '  -5.69%  ' }
    @{ Cell = 'D46'; Value = '38.75' }
    @{ Cell = 'E46'; Value = '  -2.18%  ' }
    @{ Cell = 'D47'; Value = '141.24' }
    @{ Cell = 'E47'; Value = '  -2.35%  ' }
    @{ Cell = 'E48'; Value = '  -5.60%  ' }
    @{ Cell = 'E49'; Value = '  -5.75%  ' }
    @{ Cell = 'E50'; Value = '  -5.46%  ' }
    @{ Cell = 'E51'; Value = '  -8.90%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
}
